$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") computed from regenerated save_data
$gValues = @{
    2 = 1
    3 = 0
    4 = 1
    5 = 2
    6 = 2
    7 = 3
    8 = 2
    9 = 1
    10 = 2
    11 = 2
    12 = 0
    13 = 2
    14 = 3
    15 = 3
    16 = 0
    17 = 3
    18 = 2
    19 = 0
    20 = 0
    21 = 1
    22 = 2
    23 = 1
    24 = 0
    25 = 2
    26 = 0
    27 = 3
    28 = 2
    29 = 2
    30 = 1
    31 = 1
    32 = 0
    33 = 0
    34 = 1
    35 = 0
    36 = 1
    37 = 2
    38 = 1
    39 = 1
    40 = 0
    41 = 2
    42 = 1
    43 = 1
    44 = 0
    45 = 2
    46 = 1
    47 = 1
    48 = 1
    49 = 2
    50 = 1
    51 = 3
    52 = 0
    53 = 0
    54 = 1
    55 = 2
    56 = 1
    57 = 1
    58 = 0
    59 = 2
    60 = 0
    61 = 1
    62 = 2
    63 = 1
    64 = 1
    65 = 2
    66 = 0
    67 = 2
    68 = 1
    69 = 0
    70 = 1
    71 = 2
    72 = 1
    73 = 1
    74 = 1
    75 = 1
    76 = 3
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $gValues[$row]
}
